$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 21:05"

# --- Simple numeric refreshes (no country re-sort involved) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1608420
$ws.Range("C4").Value = 15697
$ws.Range("D4").Value = 373225
$ws.Range("E4").Value = 1139433
$ws.Range("G4").Value = 826
$ws.Range("H4").Value = 95762

# Row 10: Francia
$ws.Range("B10").Value = 181826
$ws.Range("C10").Value = 251
$ws.Range("E10").Value = 89753

# Row 14: India
$ws.Range("B14").Value = 118226
$ws.Range("C14").Value = 6198
$ws.Range("D14").Value = 48553
$ws.Range("E14").Value = 66089

# Row 105: Sri Lanka
$ws.Range("B105").Value = 1048
$ws.Range("C105").Value = 20
$ws.Range("E105").Value = 435

# --- Rows 132-137: Congo's case count overtook Nepal's, so Congo moved
# up in the (descending, by total cases) ranking. Nepal, Reunion, Taiwan
# and Republica de Africa Central each slide down one row (their data is
# unchanged), Congo gets fresh totals at row 132, and Estado de Palestina
# (row 137) keeps its place but also gets fresh totals. ---

$ws.Range("A132").Value = "Congo"
$ws.Range("B132").Value = 469
$ws.Range("C132").Value = 49
$ws.Range("D132").Value = 137
$ws.Range("E132").Value = 316
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 16

$ws.Range("A133").Value = "Nepal"
$ws.Range("B133").Value = 457
$ws.Range("C133").Value = 30
$ws.Range("D133").Value = 49
$ws.Range("E133").Value = 405
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 3

$ws.Range("A134").Value = "Reunion"
$ws.Range("B134").Value = 449
$ws.Range("C134").Value = 2
$ws.Range("D134").Value = 411
$ws.Range("E134").Value = 37
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 1

$ws.Range("A135").Value = "Taiwan"
$ws.Range("B135").Value = 441
$ws.Range("C135").Value = 1
$ws.Range("D135").Value = 407
$ws.Range("E135").Value = 27
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 7

$ws.Range("A136").Value = "Republica de Africa Central"
$ws.Range("B136").Value = 436
$ws.Range("C136").Value = 18
$ws.Range("D136").Value = 18
$ws.Range("E136").Value = 418
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

$ws.Range("A137").Value = "Estado de Palestina"
$ws.Range("B137").Value = 423
$ws.Range("C137").Value = 25
$ws.Range("D137").Value = 346
$ws.Range("E137").Value = 75
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 2
